$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format to prevent Excel auto-converting price strings to numbers
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.020.44"
$ws.Range("E2").Value = "  -3.83%  "
$ws.Range("D3").Value = "3.415.06"
$ws.Range("E3").Value = "  -5.31%  "
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").Value = "568.22"
$ws.Range("E5").Value = "  -5.65%  "
$ws.Range("D6").Value = "186.23"
$ws.Range("E6").Value = "  -8.07%  "
$ws.Range("D7").Value = "0.606"
$ws.Range("E7").Value = "  -3.60%  "
$ws.Range("D8").Value = "3.409.71"
$ws.Range("E8").Value = "  -5.13%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "0.197"
$ws.Range("E10").Value = "  -8.38%  "
$ws.Range("D11").Value = "0.605"
$ws.Range("E11").Value = "  -6.46%  "
$ws.Range("D12").Value = "48.85"
$ws.Range("E12").Value = "  -9.28%  "
$ws.Range("D13").Value = "0.0000274"
$ws.Range("E13").Value = "  -9.50%  "
$ws.Range("D14").Value = "4.002.65"
$ws.Range("E14").Value = "  -4.21%  "
$ws.Range("D15").Value = "8.86"
$ws.Range("E15").Value = "  -7.85%  "
$ws.Range("D16").Value = "626.91"
$ws.Range("E16").Value = "  -7.30%  "
$ws.Range("D17").Value = "68.202.54"
$ws.Range("E17").Value = "  -3.69%  "
$ws.Range("D18").Value = "3.446.85"
$ws.Range("E18").Value = "  -3.82%  "
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("D20").Value = "17.91"
$ws.Range("E20").Value = "  -6.57%  "
$ws.Range("D21").Value = "11.89"
$ws.Range("E21").Value = "  -7.17%  "
$ws.Range("D22").Value = "0.926"
$ws.Range("E22").Value = "  -7.26%  "
$ws.Range("D23").Value = "17.74"
$ws.Range("E23").Value = "  -5.00%  "
$ws.Range("D24").Value = "5.13"
$ws.Range("E24").Value = "  -4.24%  "
$ws.Range("D25").Value = "96.06"
$ws.Range("E25").Value = "  -12.61%  "
$ws.Range("D26").Value = "4.19"
$ws.Range("E26").Value = "  -9.44%  "
$ws.Range("D27").Value = "2.78"
$ws.Range("E27").Value = "  -8.46%  "
$ws.Range("D28").Value = "9.66"
$ws.Range("E28").Value = "  -8.74%  "
$ws.Range("D29").Value = "9.01"
$ws.Range("E29").Value = "  -11.23%  "
$ws.Range("D30").Value = "31.34"
$ws.Range("E30").Value = "  -8.70%  "
$ws.Range("D31").Value = "3.97"
$ws.Range("E31").Value = "  -11.04%  "
$ws.Range("D32").Value = "6.45"
$ws.Range("E32").Value = "  -10.11%  "
$ws.Range("D33").Value = "11.34"
$ws.Range("E33").Value = "  -7.59%  "
$ws.Range("D34").Value = "567.10"
$ws.Range("E34").Value = "  +10.44%  "
$ws.Range("E35").Value = "  -7.36%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.835.13"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "58.71"
$ws.Range("E37").Value = "  -7.65%  "
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0759"
$ws.Range("E39").Value = "  -10.97%  "
$ws.Range("B40").Value = "CoreDAO"
$ws.Range("C40").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D40").Value = "3.60"
$ws.Range("E40").Value = "  +29.31%  "
$ws.Range("D41").Value = "3.43"
$ws.Range("E41").Value = "  -4.92%  "
$ws.Range("D42").Value = "2.79"
$ws.Range("E42").Value = "  -7.75%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.131"
$ws.Range("E43").Value = "  -6.65%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.359"
$ws.Range("E44").Value = "  -6.63%  "
$ws.Range("D45").Value = "32.91"
$ws.Range("E45").Value = "  -11.12%  "
$ws.Range("D46").Value = "0.0430"
$ws.Range("E46").Value = "  -8.06%  "
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("D48").Value = "2.76"
$ws.Range("E48").Value = "  -9.58%  "
$ws.Range("D49").Value = "0.133"
$ws.Range("E49").Value = "  -6.65%  "
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").Value = "7.94"
$ws.Range("E51").Value = "  -7.97%  "
